$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (target cluster "ECs" as target removed from last
# sender block; table shrinks from 20 data rows to 15 data rows, 5 senders x 3 targets)
$ws.Rows("17:21").Delete()

$data = @(
  @("ECs", "Psap", "Gpr37l1", "FAPs", 3, 1, 128.284162, 384.852486, 0.02759562102610217, 0.02759562102610216, 1, 0.3333333333333333, 0.1717676666666667, 0.515303, 0.6917552331052554, 0.6917552331052554, 22.03507117702867, 198.315640593258, 0.01908941525559559, 0.01908941525559559),
  @("ECs", "Psap", "Gpr37l1", "MuSCs", 3, 1, 128.284162, 384.852486, 0.02759562102610217, 0.02759562102610216, 1, 0.3333333333333333, 0.038916, 0.116748, 0.1567253440297696, 0.1567253440297696, 4.992306448392, 44.930758035528, 0.004324933199031006, 0.004324933199031006),
  @("ECs", "Psap", "Gpr37l1", "Resolving-Mac", 3, 1, 128.284162, 384.852486, 0.02759562102610217, 0.02759562102610216, 1, 0.3333333333333333, 0.03762333333333334, 0.11287, 0.1515194228649749, 0.1515194228649749, 4.826477788313334, 43.43830009482, 0.004181272571475568, 0.004181272571475567),
  @("FAPs", "Psap", "Gpr37l1", "FAPs", 3, 1, 166.2925923333333, 498.877777, 0.03577173741430972, 0.03577173741430972, 1, 0.3333333333333333, 0.1717676666666667, 0.515303, 0.6917552331052554, 0.6917552331052554, 28.56369056904789, 257.073215121431, 0.0247452865536158, 0.0247452865536158),
  @("FAPs", "Psap", "Gpr37l1", "MuSCs", 3, 1, 166.2925923333333, 498.877777, 0.03577173741430972, 0.03577173741430972, 1, 0.3333333333333333, 0.038916, 0.116748, 0.1567253440297696, 0.1567253440297696, 6.471442523244001, 58.24298270919601, 0.005606337852800271, 0.005606337852800272),
  @("FAPs", "Psap", "Gpr37l1", "Resolving-Mac", 3, 1, 166.2925923333333, 498.877777, 0.03577173741430972, 0.03577173741430972, 1, 0.3333333333333333, 0.03762333333333334, 0.11287, 0.1515194228649749, 0.1515194228649749, 6.256481632221112, 56.30833468999001, 0.00542011300789364, 0.00542011300789364),
  @("Inflammatory-Mac", "Psap", "Gpr37l1", "FAPs", 3, 1, 1925.076375666667, 5775.229127000001, 0.4141094058766164, 0.4141094058766165, 1, 0.3333333333333333, 0.1717676666666667, 0.515303, 0.6917552331052554, 0.6917552331052554, 330.6658772033867, 2975.992894830481, 0.2864623485932576, 0.2864623485932577),
  @("Inflammatory-Mac", "Psap", "Gpr37l1", "MuSCs", 3, 1, 1925.076375666667, 5775.229127000001, 0.4141094058766164, 0.4141094058766165, 1, 0.3333333333333333, 0.038916, 0.116748, 0.1567253440297696, 0.1567253440297696, 74.916272235444, 674.2464501189961, 0.0649014391019762, 0.06490143910197622),
  @("Inflammatory-Mac", "Psap", "Gpr37l1", "Resolving-Mac", 3, 1, 1925.076375666667, 5775.229127000001, 0.4141094058766164, 0.4141094058766165, 1, 0.3333333333333333, 0.03762333333333334, 0.11287, 0.1515194228649749, 0.1515194228649749, 72.42779017383224, 651.85011156449, 0.06274561818138258, 0.06274561818138259),
  @("MuSCs", "Psap", "Gpr37l1", "FAPs", 3, 1, 37.52106933333334, 112.563208, 0.008071278587116393, 0.008071278587116393, 1, 0.3333333333333333, 0.1717676666666667, 0.515303, 0.6917552331052554, 0.6917552331052554, 6.444906530224889, 58.004158772024, 0.005583349200488157, 0.005583349200488157),
  @("MuSCs", "Psap", "Gpr37l1", "MuSCs", 3, 1, 37.52106933333334, 112.563208, 0.008071278587116393, 0.008071278587116393, 1, 0.3333333333333333, 0.038916, 0.116748, 0.1567253440297696, 0.1567253440297696, 1.460169934176, 13.141529407584, 0.001264973913325929, 0.00126497391332593),
  @("MuSCs", "Psap", "Gpr37l1", "Resolving-Mac", 3, 1, 37.52106933333334, 112.563208, 0.008071278587116393, 0.008071278587116393, 1, 0.3333333333333333, 0.03762333333333334, 0.11287, 0.1515194228649749, 0.1515194228649749, 1.411667698551111, 12.70500928696, 0.001222955473302306, 0.001222955473302306),
  @("Resolving-Mac", "Psap", "Gpr37l1", "FAPs", 3, 1, 2391.540242666666, 7174.620728, 0.5144519570958551, 0.5144519570958552, 1, 0.3333333333333333, 0.1717676666666667, 0.515303, 0.6917552331052554, 0.6917552331052554, 410.789287222287, 3697.103585000584, 0.3558748335022981, 0.3558748335022982),
  @("Resolving-Mac", "Psap", "Gpr37l1", "MuSCs", 3, 1, 2391.540242666666, 7174.620728, 0.5144519570958551, 0.5144519570958552, 1, 0.3333333333333333, 0.038916, 0.116748, 0.1567253440297696, 0.1567253440297696, 93.06918008361599, 837.622620752544, 0.08062765996263616, 0.0806276599626362),
  @("Resolving-Mac", "Psap", "Gpr37l1", "Resolving-Mac", 3, 1, 2391.540242666666, 7174.620728, 0.5144519570958551, 0.5144519570958552, 1, 0.3333333333333333, 0.03762333333333334, 0.11287, 0.1515194228649749, 0.1515194228649749, 89.97771572992889, 809.7994415693599, 0.07794946363092081, 0.07794946363092084)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value2 = $row[$j]
  }
}
